# Word COM-interop script implementing the AC5 "Regra de Negocio" updates:
#   1) Fix the RN003 bullet text (it was split across two runs, "RN00" and
#      "3 - Os deliveries sao realizados ...", and reads correctly once merged).
#   2) Add a new RN004 bullet describing the store-warranty rule.
#   3) Add a new RN005 bullet describing the no-exchange-if-damaged rule.

$d = $word.ActiveDocument

# --- 1) Normalize "RN00" + "3 – Os deliveries ..." (two runs) into a single
#        run reading "RN003 – Os deliveries ...". A Find/Replace spanning the
#        run boundary collapses the paragraph's runs into one.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("RN003 – Os deliveries", $true, $false, $false, $false, $false, $true, 1, $false, `
    "RN003 – Os deliveries", 2) | Out-Null

# --- 2) Append the RN004 bullet after the last paragraph (same list style).
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$rn004 = $d.Paragraphs.Add($endRange)
$rn004.Range.Text = "RN004 – O cliente tem a garantia com a loja de até um mês em relação aos produtos perecíveis (alimentos e medicações) e três meses sobre em relação aos produtos não perecíveis. Caso haja expiração da validade da garantia, a loja não poderá trocar o produto."

# --- 3) Append the RN005 bullet after RN004 (same list style).
$endRange2 = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$rn005 = $d.Paragraphs.Add($endRange2)
$rn005.Range.Text = "RN005 – Produtos que foram montados, instalados ou usados de maneira incorreta no qual acarretou a danificação do produto não poderão ser trocados. "
